$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force text format on cells whose new values look like plain numbers,
# so Excel keeps them as text (matching original inlineStr cells) instead
# of converting them to numeric cells.
$textCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "46.712.32"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "2.265.23"
$ws.Range("E3").Value = "  -2.04%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "299.48"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").Value = "99.96"
$ws.Range("E6").Value = "  +2.56%  "
$ws.Range("D7").Value = "0.561"
$ws.Range("E7").Value = "  -1.85%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "0.507"
$ws.Range("E9").Value = "  -2.91%  "
$ws.Range("D10").Value = "35.03"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").Value = "0.0802"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "7.06"
$ws.Range("E12").Value = "  -3.80%  "
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D14").Value = "2.610.40"
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "2.262.85"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("D17").Value = "46.724.75"
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("D18").Value = "0.792"
$ws.Range("E18").Value = "  -3.70%  "
$ws.Range("D19").Value = "12.68"
$ws.Range("E19").Value = "  -4.26%  "
$ws.Range("E20").Value = "  +3.14%  "
$ws.Range("D21").Value = "5.80"
$ws.Range("E21").Value = "  -5.62%  "
$ws.Range("D22").Value = "65.55"
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("D23").Value = "248.03"
$ws.Range("E23").Value = "  +2.46%  "
$ws.Range("D24").Value = "2.79"
$ws.Range("E24").Value = "  -4.55%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("D26").Value = "1.85"
$ws.Range("E26").Value = "  -4.60%  "
$ws.Range("D27").Value = "41.54"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  -3.67%  "
$ws.Range("D29").Value = "9.61"
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("D30").Value = "20.23"
$ws.Range("E30").Value = "  +1.52%  "
$ws.Range("D31").Value = "2.82"
$ws.Range("E31").Value = "  +8.43%  "
$ws.Range("D32").Value = "147.02"
$ws.Range("E32").Value = "  -3.16%  "
$ws.Range("E33").Value = "  +10.06%  "
$ws.Range("D34").Value = "5.34"
$ws.Range("E34").Value = "  -5.33%  "
$ws.Range("D35").Value = "0.0767"
$ws.Range("E35").Value = "  -4.27%  "
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").Value = "  +7.85%  "
$ws.Range("D37").Value = "0.114"
$ws.Range("E37").Value = "  -2.85%  "
$ws.Range("E38").Value = "  +13.64%  "
$ws.Range("D39").Value = "1.68"
$ws.Range("E39").Value = "  -6.77%  "
$ws.Range("D40").Value = "3.84"
$ws.Range("E40").Value = "  -4.47%  "
$ws.Range("D41").Value = "0.0296"
$ws.Range("E41").Value = "  -5.91%  "
$ws.Range("D42").Value = "3.10"
$ws.Range("E42").Value = "  -6.11%  "
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("B44").Value = "BitcoinSV"
$ws.Range("C44").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D44").Value = "90.90"
$ws.Range("E44").Value = "  +16.31%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.778.04"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "1.89"
$ws.Range("E46").Value = "  -4.71%  "
$ws.Range("D47").Value = "70.96"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").Value = "0.183"
$ws.Range("E48").Value = "  -6.32%  "
$ws.Range("D49").Value = "4.79"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").Value = "93.87"
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("E51").Value = "  -2.07%  "
